$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 170 (id=168)
$ws.Range("B170").Value = 6937266
$ws.Range("G170").Value = 3
$ws.Range("H170").Value = 1
$ws.Range("I170").Value = 1
$ws.Range("J170").Value = 1
$ws.Range("L170").Value = 2.3
$ws.Range("M170").Value = 3.2
$ws.Range("N170").Value = 3.1
$ws.Range("O170").Value = 2.2
$ws.Range("P170").Value = 3.3
$ws.Range("Q170").Value = 3.3
$ws.Range("R170").Value = -0.25
$ws.Range("S170").Value = 1.925
$ws.Range("T170").Value = 1.925
$ws.Range("U170").Value = 2.5
$ws.Range("V170").Value = 2.025
$ws.Range("W170").Value = 1.825
$ws.Range("X170").Value = 1.2
$ws.Range("Y170").Value = -1
$ws.Range("Z170").Value = -1
$ws.Range("AA170").Value = 0.925
$ws.Range("AB170").Value = -1
$ws.Range("AC170").Value = 1.025
$ws.Range("AD170").Value = -1
$ws.Range("E170").Value = "Atromitos Athinon"
$ws.Range("F170").Value = "Lamia"
$ws.Range("K170").Value = "H"

# Row 171 (id=169)
$ws.Range("B171").Value = 6937268
$ws.Range("G171").Value = 1
$ws.Range("H171").Value = 2
$ws.Range("I171").Value = 1
$ws.Range("J171").Value = 0
$ws.Range("L171").Value = 8
$ws.Range("M171").Value = 5
$ws.Range("N171").Value = 1.363
$ws.Range("O171").Value = 8.5
$ws.Range("P171").Value = 5
$ws.Range("Q171").Value = 1.363
$ws.Range("R171").Value = 1.25
$ws.Range("S171").Value = 2.025
$ws.Range("T171").Value = 1.825
$ws.Range("U171").Value = 2.75
$ws.Range("V171").Value = 1.85
$ws.Range("W171").Value = 2
$ws.Range("X171").Value = -1
$ws.Range("Y171").Value = -1
$ws.Range("Z171").Value = 0.363
$ws.Range("AA171").Value = 0.5125
$ws.Range("AB171").Value = -0.5
$ws.Range("AC171").Value = 0.425
$ws.Range("AD171").Value = -0.5
$ws.Range("E171").Value = "Panetolikos"
$ws.Range("F171").Value = "Olympiakos"
$ws.Range("K171").Value = "A"

# Row 175 (id=173)
$ws.Range("B175").Value = 6937271
$ws.Range("G175").Value = 1
$ws.Range("H175").Value = 1
$ws.Range("I175").Value = 1
$ws.Range("J175").Value = 0
$ws.Range("L175").Value = 2.45
$ws.Range("M175").Value = 3.1
$ws.Range("N175").Value = 3.1
$ws.Range("O175").Value = 2
$ws.Range("P175").Value = 3.3
$ws.Range("Q175").Value = 4
$ws.Range("R175").Value = -0.5
$ws.Range("S175").Value = 2.025
$ws.Range("T175").Value = 1.825
$ws.Range("U175").Value = 2.25
$ws.Range("V175").Value = 1.85
$ws.Range("W175").Value = 2
$ws.Range("X175").Value = -1
$ws.Range("Y175").Value = 2.3
$ws.Range("Z175").Value = -1
$ws.Range("AA175").Value = -1
$ws.Range("AB175").Value = 0.825
$ws.Range("AC175").Value = -0.5
$ws.Range("AD175").Value = 0.5
$ws.Range("E175").Value = "Giannina"
$ws.Range("F175").Value = "Atromitos Athinon"
$ws.Range("K175").Value = "D"

# Row 176 (id=174)
$ws.Range("B176").Value = 6935700
$ws.Range("G176").Value = 2
$ws.Range("H176").Value = 1
$ws.Range("I176").Value = 1
$ws.Range("J176").Value = 1
$ws.Range("L176").Value = 2.6
$ws.Range("M176").Value = 3.2
$ws.Range("N176").Value = 2.875
$ws.Range("O176").Value = 2.25
$ws.Range("P176").Value = 3.3
$ws.Range("Q176").Value = 3.3
$ws.Range("R176").Value = -0.25
$ws.Range("S176").Value = 1.925
$ws.Range("T176").Value = 1.925
$ws.Range("U176").Value = 2.25
$ws.Range("V176").Value = 2
$ws.Range("W176").Value = 1.85
$ws.Range("X176").Value = 1.25
$ws.Range("Y176").Value = -1
$ws.Range("Z176").Value = -1
$ws.Range("AA176").Value = 0.925
$ws.Range("AB176").Value = -1
$ws.Range("AC176").Value = 1
$ws.Range("AD176").Value = -1
$ws.Range("E176").Value = "Panserraikos"
$ws.Range("F176").Value = "Asteras Tripolis"
$ws.Range("K176").Value = "H"

# Row 177 (id=175)
$ws.Range("B177").Value = 6935701
$ws.Range("G177").Value = 2
$ws.Range("H177").Value = 2
$ws.Range("I177").Value = 1
$ws.Range("J177").Value = 0
$ws.Range("L177").Value = 2.45
$ws.Range("M177").Value = 3.25
$ws.Range("N177").Value = 3
$ws.Range("O177").Value = 2.05
$ws.Range("P177").Value = 3.3
$ws.Range("Q177").Value = 3.8
$ws.Range("R177").Value = -0.5
$ws.Range("S177").Value = 2.05
$ws.Range("T177").Value = 1.8
$ws.Range("U177").Value = 2.25
$ws.Range("V177").Value = 1.8
$ws.Range("W177").Value = 2.05
$ws.Range("X177").Value = -1
$ws.Range("Y177").Value = 2.3
$ws.Range("Z177").Value = -1
$ws.Range("AA177").Value = -1
$ws.Range("AB177").Value = 0.8
$ws.Range("AC177").Value = 0.8
$ws.Range("AD177").Value = -1
$ws.Range("E177").Value = "Kifisias FC"
$ws.Range("F177").Value = "Panetolikos"
$ws.Range("K177").Value = "D"

# Row 178 (id=176)
$ws.Range("B178").Value = 6936863
$ws.Range("G178").Value = 2
$ws.Range("H178").Value = 2
$ws.Range("I178").Value = 1
$ws.Range("J178").Value = 1
$ws.Range("L178").Value = 8
$ws.Range("M178").Value = 4.75
$ws.Range("N178").Value = 1.4
$ws.Range("O178").Value = 5.5
$ws.Range("P178").Value = 4.75
$ws.Range("Q178").Value = 1.55
$ws.Range("R178").Value = 1
$ws.Range("S178").Value = 1.95
$ws.Range("T178").Value = 1.9
$ws.Range("U178").Value = 2.5
$ws.Range("V178").Value = 1.85
$ws.Range("W178").Value = 2
$ws.Range("X178").Value = -1
$ws.Range("Y178").Value = 3.75
$ws.Range("Z178").Value = -1
$ws.Range("AA178").Value = 0.95
$ws.Range("AB178").Value = -1
$ws.Range("AC178").Value = 0.8500000000000001
$ws.Range("AD178").Value = -1
$ws.Range("E178").Value = "OFI Crete"
$ws.Range("F178").Value = "Panathinaikos"
$ws.Range("K178").Value = "D"

# Row 179 (id=177)
$ws.Range("B179").Value = 6937269
$ws.Range("G179").Value = 3
$ws.Range("H179").Value = 3
$ws.Range("I179").Value = 1
$ws.Range("J179").Value = 1
$ws.Range("L179").Value = 4.75
$ws.Range("M179").Value = 3.75
$ws.Range("N179").Value = 1.75
$ws.Range("O179").Value = 6.5
$ws.Range("P179").Value = 4.2
$ws.Range("Q179").Value = 1.5
$ws.Range("R179").Value = 1
$ws.Range("S179").Value = 2.05
$ws.Range("T179").Value = 1.8
$ws.Range("U179").Value = 2.5
$ws.Range("V179").Value = 1.975
$ws.Range("W179").Value = 1.875
$ws.Range("X179").Value = -1
$ws.Range("Y179").Value = 3.2
$ws.Range("Z179").Value = -1
$ws.Range("AA179").Value = 1.05
$ws.Range("AB179").Value = -1
$ws.Range("AC179").Value = 0.9750000000000001
$ws.Range("AD179").Value = -1
$ws.Range("E179").Value = "Aris Salonika"
$ws.Range("F179").Value = "AEK Athens"
$ws.Range("K179").Value = "D"

# Row 180 (id=178)
$ws.Range("B180").Value = 6937270
$ws.Range("G180").Value = 3
$ws.Range("H180").Value = 0
$ws.Range("I180").Value = 2
$ws.Range("J180").Value = 0
$ws.Range("L180").Value = 1.125
$ws.Range("M180").Value = 9
$ws.Range("N180").Value = 19
$ws.Range("O180").Value = 1.111
$ws.Range("P180").Value = 9
$ws.Range("Q180").Value = 21
$ws.Range("R180").Value = -2.25
$ws.Range("S180").Value = 1.875
$ws.Range("T180").Value = 1.975
$ws.Range("U180").Value = 3.25
$ws.Range("V180").Value = 2
$ws.Range("W180").Value = 1.85
$ws.Range("X180").Value = 0.111
$ws.Range("Y180").Value = -1
$ws.Range("Z180").Value = -1
$ws.Range("AA180").Value = 0.875
$ws.Range("AB180").Value = -1
$ws.Range("AC180").Value = -0.5
$ws.Range("AD180").Value = 0.425
$ws.Range("E180").Value = "Olympiakos"
$ws.Range("F180").Value = "Volos NFC"
$ws.Range("K180").Value = "H"
